$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -0.03842403594624615
$ws.Range("B2").Value = 0.01336580053704416
$ws.Range("C2").Value = -0.02111056220124395
$ws.Range("D2").Value = 0.00503848044508
$ws.Range("E2").Value = -0.008560712252869118
$ws.Range("F2").Value = "Nguyễn Ngọc Lâm"

$ws.Range("A3").Value = -0.04110429668360206
$ws.Range("B3").Value = 0.05035112084275228
$ws.Range("C3").Value = -0.01174729989530701
$ws.Range("D3").Value = 0.004176293938481895
$ws.Range("E3").Value = -0.001450596453483963
$ws.Range("F3").Value = "Phạm Thị Hòa"

$ws.Range("A4").Value = -0.02920324986108187
$ws.Range("B4").Value = 0.01246068639759188
$ws.Range("C4").Value = 0.01645955398744802
$ws.Range("D4").Value = 0.004946250690690694
$ws.Range("E4").Value = 0.008666237927310878
$ws.Range("F4").Value = "Nguyễn Minh Khôi"

$ws.Range("A5").Value = -0.02480676646845765
$ws.Range("B5").Value = 0.03087050028332895
$ws.Range("C5").Value = 0.02525817687221753
$ws.Range("D5").Value = -0.02179462479754733
$ws.Range("E5").Value = -0.003652766469384778
$ws.Range("F5").Value = "Trần Thị Thùy Linh"

$ws.Range("A6").Value = 0.3934377193377001
$ws.Range("B6").Value = -0.0002005818893808036
$ws.Range("C6").Value = -0.00238421476150934
$ws.Range("D6").Value = 0.0001985690767570844
$ws.Range("E6").Value = 0.0001089655826942301
$ws.Range("F6").Value = "Nguyễn Thị Ngọc Thoa"

$ws.Range("A7").Value = -0.02967628178700396
$ws.Range("B7").Value = -0.01603961925517172
$ws.Range("C7").Value = 0.01802280440629499
$ws.Range("D7").Value = 0.008191921894608059
$ws.Range("E7").Value = 0.02120467639060371
$ws.Range("F7").Value = "Lê Xuân Quý"

$ws.Range("A8").Value = -0.04139033625284454
$ws.Range("B8").Value = -0.004370600022166111
$ws.Range("C8").Value = 0.0009907102970992405
$ws.Range("D8").Value = -0.0006710254545786845
$ws.Range("E8").Value = 0.001174806176446212
$ws.Range("F8").Value = "Ma Chí Định"

$ws.Range("A9").Value = -0.04189526835967602
$ws.Range("B9").Value = -0.005724881702870478
$ws.Range("C9").Value = -0.01685239313026183
$ws.Range("D9").Value = 0.0004961700871757662
$ws.Range("E9").Value = 0.004659207225906152
$ws.Range("F9").Value = "Ma Chí Định"

$ws.Range("A10").Value = -0.03287580758557897
$ws.Range("B10").Value = -0.02420556564215289
$ws.Range("C10").Value = -0.001501964982323434
$ws.Range("D10").Value = -0.007763312207248805
$ws.Range("E10").Value = -0.002144893170773051
$ws.Range("F10").Value = "Trần Đức Phụng"

$ws.Range("A11").Value = -0.04584315835868628
$ws.Range("B11").Value = -0.01030489015215447
$ws.Range("C11").Value = -0.02235134079560485
$ws.Range("D11").Value = 0.00330848918898061
$ws.Range("E11").Value = 0.01006939032269512
$ws.Range("F11").Value = "Phạm Thị Hương"

$ws.Range("A12").Value = -0.03096259279511591
$ws.Range("B12").Value = -0.01500062937872642
$ws.Range("C12").Value = 0.01810810838872735
$ws.Range("D12").Value = 0.02665995140111898
$ws.Range("E12").Value = -0.01989623417779658
$ws.Range("F12").Value = "Phạm Thị Hương"

$ws.Range("A13").Value = -0.03725592523940672
$ws.Range("B13").Value = -0.03120134001809426
$ws.Range("C13").Value = -0.00289157818553679
$ws.Range("D13").Value = -0.02278716426351822
$ws.Range("E13").Value = -0.01017808110134885
$ws.Range("F13").Value = "Trần Quốc Việt"
